$p = $ppt.ActivePresentation
$s = $p.Slides.Item(18)
$sh = $s.Shapes.Item(3)
$tr = $sh.TextFrame.TextRange

# Merge "computes " + "its activations" -> "computes its activations"
# (chars 4..27, length 24) picks up the rPr (F0FF33, dirty="0") of the
# first sub-run ("computes ").
$cA = $tr.Characters(4, 24)
$cA.Text = "computes its activations"

# Merge " " + "just like the " -> " just like the "
# Insert a leading space onto the "just like the " run (which already
# carries dirty="0") so the merged run keeps that rPr, then delete the
# now-redundant standalone space run that preceded it.
$cJustLikeThe = $tr.Characters(29, 14)
[void]$cJustLikeThe.InsertBefore(" ")
$cOldSpace = $tr.Characters(28, 1)
$cOldSpace.Text = ""
